# Revert "Powerpoint writer: consolidate text run nodes."
#
# Before this change the deck's writer emitted one <a:r> run per
# whitespace-delimited token *plus* a separate run for the single
# space that followed it (e.g. "A" / " " / "slide"). A later commit
# started consolidating those into fewer, merged runs (e.g. "A ").
# This script reverts that consolidation on the two paragraphs that
# were affected, splitting each merged run back into
# "<word>" + " " runs without altering the visible text.
#
# Mechanism notes (figured out empirically against this COM host):
#   - TextRange.Characters(start, len) returns a sub-range.
#   - Calling .InsertAfter(text) on a sub-range that spans from the
#     start of a run through an *interior* character inserts the new
#     text right after that sub-range, splitting the run there; the
#     newly inserted text ends up as its own trailing run (it does
#     not get silently re-merged with its neighbours).
#   - Deleting that placeholder afterwards (by setting the 1-char
#     sub-range's .Text = "") removes it again while leaving the
#     split in place, because deleting interior run text doesn't
#     trigger the "merge adjacent same-format runs" pass that
#     assigning the *whole* TextRange's .Text does.
#   - So: insert a 1-character placeholder immediately after the
#     prefix we want isolated, then delete that placeholder. Net
#     text is unchanged, but the run boundary remains.

function Split-RunAt {
    # Splits the run(s) under $tr so that the first $prefixLen
    # characters starting at 1-based $runStart become their own run.
    param($tr, [int]$runStart, [int]$prefixLen)

    $prefix = $tr.Characters($runStart, $prefixLen)
    $prefix.InsertAfter("X") | Out-Null

    $placeholderPos = $runStart + $prefixLen
    $placeholder = $tr.Characters($placeholderPos, 1)
    $placeholder.Text = ""
}

function Split-IntoWordsAndSpaces {
    # Given a TextRange whose text is a run of words separated by
    # single spaces (e.g. "Just an image"), split it so each word and
    # each separating space becomes its own run, left to right.
    param($tr)

    $text = $tr.Text
    $words = $text.Split(" ")

    $pos = 1
    for ($i = 0; $i -lt $words.Length - 1; $i++) {
        $wordLen = $words[$i].Length
        Split-RunAt $tr $pos $wordLen   # splits off the word
        $pos += $wordLen
        Split-RunAt $tr $pos 1          # splits off the following space
        $pos += 1
    }
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if (-not $shp.HasTextFrame) { continue }

    $tr = $shp.TextFrame.TextRange
    if ($shp.Name -eq "Title 1" -or $shp.Name -eq "TextBox 3") {
        Split-IntoWordsAndSpaces $tr
    }
}
